$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O header: "Internal Assignment" - bold, size 12 font, same family as the
# other header cells (K4:N4 use bold size 11 "Calibri").
$header = $ws.Cells.Item(4, 15)
$header.Value = "Internal Assignment"
$header.Font.Bold = $true
$header.Font.Size = 12
$header.Font.Name = "Calibri"
$header.Font.Color = 0

# Data rows 5-7: "Internal Assignment" column gets the same "FALSE" value as the
# "Unique" column (L) directly to its left, using the plain data-row style.
$ws.Cells.Item(5, 15).Value = "FALSE"
$ws.Cells.Item(6, 15).Value = "FALSE"
$ws.Cells.Item(7, 15).Value = "FALSE"

# Keep the active selection on the newly added bottom-right cell, matching the
# authored workbook state.
$ws.Range("O12").Select()
